$wb = $excel.ActiveWorkbook

# --- Basis sheet: mark several user stories as "Afgerond" (Done) ---
$basis = $wb.Worksheets.Item("Basis")
$basisRows = 7,8,12,17,19,20,21,24,25,26,27,28,29,30
foreach ($r in $basisRows) {
    $basis.Range("A$r").Value = "Afgerond"
}

# --- Extra sheet: mark two user stories as "Afgerond" (Done) ---
$extra = $wb.Worksheets.Item("Extra")
$extraRows = 2,12
foreach ($r in $extraRows) {
    $extra.Range("A$r").Value = "Afgerond"
}

# --- Update view/selection state to match the saved workbook ---
$basis.Select()
$basis.Range("B2:B30").Select()
$excel.ActiveWindow.Zoom = 85

$extra.Select()
$extra.Range("B12,B2").Select()
